$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 22283.658
$ws.Range("J17").Value = 22839.05
$ws.Range("L17").Value = 68517.14999999999
$ws.Range("N17").Value = -68853.14999999999
$ws.Range("H69").Value = 2906.7144
$ws.Range("I69").Value = 2725
$ws.Range("K69").Value = 8175
$ws.Range("M69").Value = -7301
$ws.Range("H72").Value = 2906.7144
$ws.Range("I72").Value = 2725
$ws.Range("K72").Value = 24525
$ws.Range("M72").Value = -20157
$ws.Range("H129").Value = 947.1111
$ws.Range("I129").Value = 816.375
$ws.Range("K129").Value = 2449.125
$ws.Range("M129").Value = 2550.875
$ws.Range("H132").Value = 3487.7805
$ws.Range("I132").Value = 1689.9231
$ws.Range("J132").Value = 6604.067
$ws.Range("K132").Value = 5069.7693
$ws.Range("L132").Value = 19812.201
$ws.Range("M132").Value = -2539.7693
$ws.Range("N132").Value = -24872.201
$ws.Range("H135").Value = 34150.867
$ws.Range("I135").Value = 774.0769
$ws.Range("J135").Value = 251100
$ws.Range("K135").Value = 6966.6921
$ws.Range("L135").Value = 2259900
$ws.Range("M135").Value = -4431.6921
$ws.Range("N135").Value = -2264970
$ws.Range("H137").Value = 26297.625
$ws.Range("I137").Value = 1297.8
$ws.Range("K137").Value = 3893.4
$ws.Range("M137").Value = -1343.4
$ws.Range("H138").Value = 2348.394
$ws.Range("J138").Value = 3158.625
$ws.Range("L138").Value = 9475.875
$ws.Range("N138").Value = -19755.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2229.6897
$ws.Range("I2").Value = 1965.9474
$ws.Range("J2").Value = 2730.8
$ws.Range("K2").Value = 1965.9474
$ws.Range("L2").Value = 2730.8
$ws.Range("M2").Value = -1852.9474
$ws.Range("N2").Value = -2956.8
$ws.Range("H32").Value = 39181.508
$ws.Range("I32").Value = 21219.82
$ws.Range("K32").Value = 21219.82
$ws.Range("M32").Value = -20932.82
$ws.Range("H46").Value = 9164.666999999999
$ws.Range("J46").Value = 8874.5
$ws.Range("L46").Value = 8874.5
$ws.Range("N46").Value = -9512.5
$ws.Range("H110").Value = 2872.3
$ws.Range("I110").Value = 2340.375
$ws.Range("K110").Value = 2340.375
$ws.Range("M110").Value = -295.375
$ws.Range("H116").Value = 2229.6897
$ws.Range("I116").Value = 1965.9474
$ws.Range("J116").Value = 2730.8
$ws.Range("K116").Value = 1965.9474
$ws.Range("L116").Value = 2730.8
$ws.Range("M116").Value = 328.0526
$ws.Range("N116").Value = -7318.8
$ws.Range("H132").Value = 35457
$ws.Range("I132").Value = 41802.133
$ws.Range("J132").Value = 3731.3333
$ws.Range("K132").Value = 125406.399
$ws.Range("L132").Value = 11193.9999
$ws.Range("M132").Value = -122876.399
$ws.Range("N132").Value = -16253.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2229.6897
$ws.Range("I3").Value = 1965.9474
$ws.Range("J3").Value = 2730.8
$ws.Range("K3").Value = 1965.9474
$ws.Range("L3").Value = 2730.8
$ws.Range("M3").Value = -1851.9474
$ws.Range("N3").Value = -2958.8
$ws.Range("H80").Value = 1236
$ws.Range("I80").Value = 888.625
$ws.Range("K80").Value = 888.625
$ws.Range("M80").Value = 109.375
$ws.Range("H83").Value = 1236
$ws.Range("I83").Value = 888.625
$ws.Range("K83").Value = 4443.125
$ws.Range("M83").Value = 548.875
$ws.Range("H107").Value = 5523.7144
$ws.Range("I107").Value = 3042.2173
$ws.Range("K107").Value = 3042.2173
$ws.Range("M107").Value = -1122.2173

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1101.5454
$ws.Range("I16").Value = 859.8570999999999
$ws.Range("K16").Value = 859.8570999999999
$ws.Range("M16").Value = -572.8570999999999
$ws.Range("H31").Value = 1610.1923
$ws.Range("I31").Value = 1707.8
$ws.Range("K31").Value = 1707.8
$ws.Range("M31").Value = -1412.8
$ws.Range("H34").Value = 1610.1923
$ws.Range("I34").Value = 1707.8
$ws.Range("K34").Value = 1707.8
$ws.Range("M34").Value = -1505.8
$ws.Range("H99").Value = 168433
$ws.Range("I99").Value = 201619.6
$ws.Range("K99").Value = 201619.6
$ws.Range("M99").Value = -200121.6
$ws.Range("H107").Value = 2159.5625
$ws.Range("I107").Value = 2270.6667
$ws.Range("J107").Value = 2092.9
$ws.Range("K107").Value = 2270.6667
$ws.Range("L107").Value = 2092.9
$ws.Range("M107").Value = -350.6667000000002
$ws.Range("N107").Value = -5932.9
$ws.Range("H113").Value = 1101.5454
$ws.Range("I113").Value = 859.8570999999999
$ws.Range("K113").Value = 859.8570999999999
$ws.Range("M113").Value = 1310.1429
$ws.Range("H126").Value = 168433
$ws.Range("I126").Value = 201619.6
$ws.Range("K126").Value = 604858.8
$ws.Range("M126").Value = -602388.8
$ws.Range("H132").Value = 3165
$ws.Range("I132").Value = 2998.182
$ws.Range("K132").Value = 8994.545999999998
$ws.Range("M132").Value = -6464.545999999998
$ws.Range("H134").Value = 3119.5386
$ws.Range("I134").Value = 3265.5
$ws.Range("K134").Value = 9796.5
$ws.Range("M134").Value = -7261.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 74.5
$ws.Range("I23").Value = 48.666668
$ws.Range("J23").Value = 100.333336
$ws.Range("K23").Value = 146.000004
$ws.Range("L23").Value = 301.000008
$ws.Range("M23").Value = 88.99999600000001
$ws.Range("N23").Value = -771.000008
$ws.Range("H29").Value = 333642
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("H46").Value = 125306.375
$ws.Range("I46").Value = 143175.86
$ws.Range("J46").Value = 220
$ws.Range("K46").Value = 429527.58
$ws.Range("L46").Value = 660
$ws.Range("M46").Value = -429436.58
$ws.Range("N46").Value = -842
$ws.Range("H62").Value = 174832.5
$ws.Range("I62").Value = 338665.34
$ws.Range("K62").Value = 1015996.02
$ws.Range("M62").Value = -1015310.02
$ws.Range("H65").Value = 174832.5
$ws.Range("I65").Value = 338665.34
$ws.Range("K65").Value = 3047988.06
$ws.Range("M65").Value = -3044556.06
$ws.Range("H120").Value = 17500
$ws.Range("I120").Value = 15000
$ws.Range("K120").Value = 45000
$ws.Range("M120").Value = -40162
$ws.Range("H131").Value = 113996.336
$ws.Range("J131").Value = 3285.5715
$ws.Range("L131").Value = 9856.7145
$ws.Range("N131").Value = -19936.7145
$ws.Range("H137").Value = 5265990.5
$ws.Range("I137").Value = 8334566
$ws.Range("J137").Value = 5575.7144
$ws.Range("K137").Value = 25003698
$ws.Range("L137").Value = 16727.1432
$ws.Range("M137").Value = -24998598
$ws.Range("N137").Value = -26927.1432
$ws.Range("N29").Value = $null

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5212.857
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 5747.5
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 5747.5
$ws.Range("M80").Value = -3502
$ws.Range("N80").Value = -7743.5
$ws.Range("H83").Value = 5212.857
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 5747.5
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 28737.5
$ws.Range("M83").Value = -17508
$ws.Range("N83").Value = -38721.5
$ws.Range("H113").Value = 29414242
$ws.Range("I113").Value = 38463824
$ws.Range("J113").Value = 3102
$ws.Range("K113").Value = 38463824
$ws.Range("L113").Value = 3102
$ws.Range("M113").Value = -38461654
$ws.Range("N113").Value = -7442
$ws.Range("H132").Value = 4069.4285
$ws.Range("I132").Value = 4472.75
$ws.Range("J132").Value = 1649.5
$ws.Range("K132").Value = 13418.25
$ws.Range("L132").Value = 4948.5
$ws.Range("M132").Value = -10888.25
$ws.Range("N132").Value = -10008.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 420.2857
$ws.Range("I9").Value = 274.375
$ws.Range("J9").Value = 614.8333
$ws.Range("K9").Value = 274.375
$ws.Range("L9").Value = 614.8333
$ws.Range("M9").Value = -50.375
$ws.Range("N9").Value = -1062.8333
$ws.Range("H46").Value = 2497.625
$ws.Range("I46").Value = 2497.625
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2497.625
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2309.625
$ws.Range("H61").Value = 650.25
$ws.Range("I61").Value = 650.25
$ws.Range("K61").Value = 650.25
$ws.Range("M61").Value = -448.25
$ws.Range("H113").Value = 650.25
$ws.Range("I113").Value = 650.25
$ws.Range("K113").Value = 650.25
$ws.Range("M113").Value = 1519.75
$ws.Range("H132").Value = 2853.8
$ws.Range("I132").Value = 2227.55
$ws.Range("J132").Value = 5358.8
$ws.Range("K132").Value = 6682.650000000001
$ws.Range("L132").Value = 16076.4
$ws.Range("M132").Value = -4152.650000000001
$ws.Range("N132").Value = -21136.4
$ws.Range("H136").Value = 3125.8
$ws.Range("I136").Value = 2241
$ws.Range("K136").Value = 6723
$ws.Range("M136").Value = -4173
$ws.Range("N46").Value = $null

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10205.7
$ws.Range("I132").Value = 13098
$ws.Range("J132").Value = 7674.9375
$ws.Range("K132").Value = 39294
$ws.Range("L132").Value = 23024.8125
$ws.Range("M132").Value = -36764
$ws.Range("N132").Value = -28084.8125
$ws.Range("H136").Value = 952.8889
$ws.Range("I136").Value = 973.64703
$ws.Range("J136").Value = 600
$ws.Range("K136").Value = 2920.94109
$ws.Range("L136").Value = 1800
$ws.Range("M136").Value = -370.9410899999998
$ws.Range("N136").Value = -6900
